$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the weekly Brocoli price rows down by two (336-370 -> 338-372),
# insert a new top entry pair (336-337, week of 2021-09-10) and
# duplicate the former last pair (369-370) as the new tail (371-372).

$ws.Cells.Item(336, 4).Value = 44449
$ws.Cells.Item(336, 10).Value = 3400
$ws.Cells.Item(336, 11).Value = 600
$ws.Cells.Item(336, 12).Value = 700
$ws.Cells.Item(336, 13).Value = 650
$ws.Cells.Item(336, 16).Value = 650
$ws.Cells.Item(337, 4).Value = 44449
$ws.Cells.Item(337, 10).Value = 1600
$ws.Cells.Item(337, 11).Value = 500
$ws.Cells.Item(337, 12).Value = 550
$ws.Cells.Item(337, 13).Value = 525
$ws.Cells.Item(337, 16).Value = 525
$ws.Cells.Item(338, 4).Value = 44161
$ws.Cells.Item(338, 10).Value = 2200
$ws.Cells.Item(338, 11).Value = 450
$ws.Cells.Item(338, 12).Value = 500
$ws.Cells.Item(338, 13).Value = 475
$ws.Cells.Item(338, 16).Value = 475
$ws.Cells.Item(339, 4).Value = 44161
$ws.Cells.Item(339, 10).Value = 1480
$ws.Cells.Item(339, 11).Value = 350
$ws.Cells.Item(339, 12).Value = 400
$ws.Cells.Item(339, 13).Value = 375
$ws.Cells.Item(339, 16).Value = 375
$ws.Cells.Item(340, 4).Value = 44438
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 3200
$ws.Cells.Item(340, 11).Value = 600
$ws.Cells.Item(340, 12).Value = 700
$ws.Cells.Item(340, 13).Value = 650
$ws.Cells.Item(340, 16).Value = 650
$ws.Cells.Item(341, 4).Value = 44438
$ws.Cells.Item(341, 9).Value = "Segunda"
$ws.Cells.Item(341, 10).Value = 1500
$ws.Cells.Item(341, 11).Value = 500
$ws.Cells.Item(341, 12).Value = 550
$ws.Cells.Item(341, 13).Value = 525
$ws.Cells.Item(341, 16).Value = 525
$ws.Cells.Item(342, 4).Value = 44251
$ws.Cells.Item(342, 10).Value = 2000
$ws.Cells.Item(342, 11).Value = 750
$ws.Cells.Item(342, 12).Value = 800
$ws.Cells.Item(342, 13).Value = 775
$ws.Cells.Item(342, 16).Value = 775
$ws.Cells.Item(343, 4).Value = 44428
$ws.Cells.Item(343, 10).Value = 3400
$ws.Cells.Item(343, 11).Value = 650
$ws.Cells.Item(343, 13).Value = 675
$ws.Cells.Item(343, 16).Value = 675
$ws.Cells.Item(344, 4).Value = 44428
$ws.Cells.Item(344, 11).Value = 550
$ws.Cells.Item(344, 12).Value = 600
$ws.Cells.Item(344, 13).Value = 575
$ws.Cells.Item(344, 16).Value = 575
$ws.Cells.Item(345, 4).Value = 44442
$ws.Cells.Item(345, 10).Value = 3200
$ws.Cells.Item(346, 4).Value = 44442
$ws.Cells.Item(346, 10).Value = 1600
$ws.Cells.Item(347, 4).Value = 44435
$ws.Cells.Item(347, 10).Value = 13960
$ws.Cells.Item(347, 11).Value = 600
$ws.Cells.Item(347, 13).Value = 650
$ws.Cells.Item(347, 16).Value = 650
$ws.Cells.Item(348, 4).Value = 44435
$ws.Cells.Item(348, 10).Value = 7500
$ws.Cells.Item(349, 4).Value = 44319
$ws.Cells.Item(349, 10).Value = 3000
$ws.Cells.Item(349, 11).Value = 650
$ws.Cells.Item(349, 12).Value = 700
$ws.Cells.Item(349, 13).Value = 675
$ws.Cells.Item(349, 16).Value = 675
$ws.Cells.Item(350, 4).Value = 44319
$ws.Cells.Item(350, 10).Value = 1440
$ws.Cells.Item(350, 11).Value = 500
$ws.Cells.Item(350, 12).Value = 550
$ws.Cells.Item(350, 13).Value = 525
$ws.Cells.Item(350, 16).Value = 525
$ws.Cells.Item(351, 4).Value = 44175
$ws.Cells.Item(351, 10).Value = 2200
$ws.Cells.Item(351, 11).Value = 550
$ws.Cells.Item(351, 12).Value = 600
$ws.Cells.Item(351, 13).Value = 575
$ws.Cells.Item(351, 16).Value = 575
$ws.Cells.Item(352, 4).Value = 44175
$ws.Cells.Item(352, 10).Value = 1500
$ws.Cells.Item(352, 11).Value = 450
$ws.Cells.Item(352, 12).Value = 500
$ws.Cells.Item(352, 13).Value = 475
$ws.Cells.Item(352, 16).Value = 475
$ws.Cells.Item(353, 4).Value = 44376
$ws.Cells.Item(353, 10).Value = 2460
$ws.Cells.Item(353, 11).Value = 600
$ws.Cells.Item(353, 12).Value = 700
$ws.Cells.Item(353, 13).Value = 650
$ws.Cells.Item(353, 16).Value = 650
$ws.Cells.Item(354, 4).Value = 44376
$ws.Cells.Item(354, 10).Value = 1400
$ws.Cells.Item(354, 11).Value = 500
$ws.Cells.Item(354, 12).Value = 550
$ws.Cells.Item(354, 13).Value = 525
$ws.Cells.Item(354, 16).Value = 525
$ws.Cells.Item(355, 4).Value = 44279
$ws.Cells.Item(355, 10).Value = 3400
$ws.Cells.Item(355, 11).Value = 700
$ws.Cells.Item(355, 12).Value = 750
$ws.Cells.Item(355, 13).Value = 725
$ws.Cells.Item(355, 16).Value = 725
$ws.Cells.Item(356, 4).Value = 44279
$ws.Cells.Item(356, 10).Value = 1800
$ws.Cells.Item(356, 11).Value = 600
$ws.Cells.Item(356, 12).Value = 650
$ws.Cells.Item(356, 13).Value = 625
$ws.Cells.Item(356, 16).Value = 625
$ws.Cells.Item(357, 4).Value = 44412
$ws.Cells.Item(358, 4).Value = 44412
$ws.Cells.Item(358, 10).Value = 1700
$ws.Cells.Item(359, 4).Value = 44223
$ws.Cells.Item(359, 10).Value = 3200
$ws.Cells.Item(360, 4).Value = 44223
$ws.Cells.Item(360, 11).Value = 550
$ws.Cells.Item(360, 12).Value = 600
$ws.Cells.Item(360, 13).Value = 575
$ws.Cells.Item(360, 16).Value = 575
$ws.Cells.Item(361, 4).Value = 44314
$ws.Cells.Item(361, 10).Value = 3400
$ws.Cells.Item(361, 11).Value = 650
$ws.Cells.Item(361, 13).Value = 675
$ws.Cells.Item(361, 16).Value = 675
$ws.Cells.Item(362, 4).Value = 44314
$ws.Cells.Item(362, 10).Value = 1800
$ws.Cells.Item(363, 4).Value = 44448
$ws.Cells.Item(363, 10).Value = 2000
$ws.Cells.Item(363, 11).Value = 600
$ws.Cells.Item(363, 12).Value = 700
$ws.Cells.Item(363, 13).Value = 650
$ws.Cells.Item(363, 16).Value = 650
$ws.Cells.Item(364, 4).Value = 44448
$ws.Cells.Item(364, 10).Value = 1300
$ws.Cells.Item(364, 11).Value = 500
$ws.Cells.Item(364, 12).Value = 550
$ws.Cells.Item(364, 13).Value = 525
$ws.Cells.Item(364, 16).Value = 525
$ws.Cells.Item(365, 4).Value = 44167
$ws.Cells.Item(365, 10).Value = 2900
$ws.Cells.Item(365, 11).Value = 450
$ws.Cells.Item(365, 12).Value = 500
$ws.Cells.Item(365, 13).Value = 475
$ws.Cells.Item(365, 16).Value = 475
$ws.Cells.Item(366, 4).Value = 44167
$ws.Cells.Item(366, 10).Value = 1600
$ws.Cells.Item(366, 11).Value = 350
$ws.Cells.Item(366, 12).Value = 400
$ws.Cells.Item(366, 13).Value = 375
$ws.Cells.Item(366, 16).Value = 375
$ws.Cells.Item(367, 4).Value = 44238
$ws.Cells.Item(367, 10).Value = 2400
$ws.Cells.Item(367, 11).Value = 650
$ws.Cells.Item(367, 13).Value = 675
$ws.Cells.Item(367, 16).Value = 675
$ws.Cells.Item(368, 4).Value = 44238
$ws.Cells.Item(368, 11).Value = 550
$ws.Cells.Item(368, 12).Value = 600
$ws.Cells.Item(368, 13).Value = 575
$ws.Cells.Item(368, 16).Value = 575
$ws.Cells.Item(369, 4).Value = 44399
$ws.Cells.Item(369, 10).Value = 2200
$ws.Cells.Item(370, 4).Value = 44399
$ws.Cells.Item(370, 10).Value = 1400
$ws.Cells.Item(371, 1).Value = 8
$ws.Cells.Item(371, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(371, 3).Value = "Coquimbo"
$ws.Cells.Item(371, 4).Value = 44400
$ws.Cells.Item(371, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(371, 5).Value = 4
$ws.Cells.Item(371, 6).Value = 100112023
$ws.Cells.Item(371, 7).Value = "Brócoli"
$ws.Cells.Item(371, 8).Value = "Sin especificar"
$ws.Cells.Item(371, 9).Value = "Primera"
$ws.Cells.Item(371, 10).Value = 3400
$ws.Cells.Item(371, 11).Value = 600
$ws.Cells.Item(371, 12).Value = 700
$ws.Cells.Item(371, 13).Value = 650
$ws.Cells.Item(371, 14).Value = "`$/unidad"
$ws.Cells.Item(371, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(371, 16).Value = 650
$ws.Cells.Item(371, 17).Value = 1
$ws.Cells.Item(371, 18).Value = "Hortaliza"
$ws.Cells.Item(372, 1).Value = 8
$ws.Cells.Item(372, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value = "Coquimbo"
$ws.Cells.Item(372, 4).Value = 44400
$ws.Cells.Item(372, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(372, 5).Value = 4
$ws.Cells.Item(372, 6).Value = 100112023
$ws.Cells.Item(372, 7).Value = "Brócoli"
$ws.Cells.Item(372, 8).Value = "Sin especificar"
$ws.Cells.Item(372, 9).Value = "Segunda"
$ws.Cells.Item(372, 10).Value = 1600
$ws.Cells.Item(372, 11).Value = 500
$ws.Cells.Item(372, 12).Value = 550
$ws.Cells.Item(372, 13).Value = 525
$ws.Cells.Item(372, 14).Value = "`$/unidad"
$ws.Cells.Item(372, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(372, 16).Value = 525
$ws.Cells.Item(372, 17).Value = 1
$ws.Cells.Item(372, 18).Value = "Hortaliza"
